$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidential disclosure date (shared string used in A59)
$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-29 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-56
$data = @{
    2 = @(0.01535762295026291, 0.01057854844866157)
    3 = @(0.05363426958282014, 0.00370391788347546)
    4 = @(0.01506753110500739, -0.04407163441745576)
    5 = @(0.00945489127634922, 0.006560818790185108)
    6 = @(0.01495655616409813, 0.01487696538109007)
    7 = @(0.0191548082746613, 0.008447488584474749)
    8 = @(0.004104273890620455, 0.003544223444773564)
    9 = @(0.006488405181233275, 0.01428325589399426)
    10 = @(0.01420699456629098, 0.008939974457215838)
    11 = @(0.008041806226649474, 0.0230638691761802)
    12 = @(0.01481320060879277, 0.02286432160804019)
    13 = @(0.003017848448985629, -0.02122302158273381)
    14 = @(0.00623668002180915, 0.00820568927789922)
    15 = @(0.01369250258199732, 0.0194442619720161)
    16 = @(0.01017616636043104, 0.01275545192703342)
    17 = @(0.02188064189171271, -0.01186451489443141)
    18 = @(0.008601519413806666, -0.04755414205664099)
    19 = @(0.01607666488852431, 0.01370539572786766)
    20 = @(0.01119984660769487, 0.009138742730545601)
    21 = @(0.007105870000579688, 0.0122739018087854)
    22 = @(0.01309615960020239, 0.01286472148541118)
    23 = @(0.0192608827011313, 0.003347826086956607)
    24 = @(0.009566672631026962, 0.01062429484768712)
    25 = @(0.02064903095583635, 0.006560962274466897)
    26 = @(0.01292051647824344, 0.02696980390179959)
    27 = @(0.02176867444154993, 0.004576443852362777)
    28 = @(0.05800348141230221, -0.0007486150621351939)
    29 = @(0.0207755758856602, 0.01780883678990075)
    30 = @(0.03076158367629956, 0.009235733010687586)
    31 = @(0.01567262057459896, 0.04473896170462033)
    32 = @(0.01324751726136425, 0.01983985765124552)
    33 = @(0.01921175349308204, 0.02359641985353944)
    34 = @(0.0439006693017167, 0.01429395008138923)
    35 = @(0.01079353813203493, 0.009310344827586192)
    36 = @(0.009720511565322632, -0.009291521486643362)
    37 = @(0.01080284290628669, -0.001291989664082727)
    38 = @(0.007195381928881564, -0.000387947756368856)
    39 = @(0.0115465735122294, 0.02697976517611789)
    40 = @(0.01680336775758631, 0.006847974955977332)
    41 = @(0.0170343742863432, 0.003390299591780277)
    42 = @(0.03442208186693968, -0.01408336486997885)
    43 = @(0.01120617385418606, 0.007231804795961194)
    44 = @(0.02244606200707764, 0.01460698222317425)
    45 = @(0.01281825700921666, 0.01543989547038338)
    46 = @(0.008022297216634963, 0.0301680643028639)
    47 = @(0.01296899435209508, -0.02079207920792092)
    48 = @(0.009660588819141334, 0.0107489597780861)
    49 = @(0.01519060225244392, 0.01352477091194193)
    50 = @(0.008468492157920754, 0.01220347425440504)
    51 = @(0.01176365389552283, 0.01663160004007613)
    52 = @(0.008785598864420537, -0.02364956700710641)
    53 = @(0.009522940192043717, 0.02056449774292091)
    54 = @(0.1347658839707677, -0.0001970055161544249)
    55 = @(0.04395854499756262, 0.004797888928871119)
    56 = @(0.9999999999999999, 0.005952818653392589)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals[0]
    $ws.Range("E$row").Value = $vals[1]
}

$ws.Protect()
